# Respondendo questao 4: append the SQL "SELECT ... FROM ... GROUP BY ... ORDER BY ..."
# code listing (with its original Consolas/syntax-highlighting run formatting)
# as a block of new paragraphs at the end of the document, right after the
# existing trailing "Resposta:" / " " paragraphs and before the sectPr.
#
# We build the new paragraphs from their exact target WordprocessingML and
# hand them to Range.InsertXML on a *collapsed* range positioned at the very
# end of the document's main story (just before the final paragraph mark),
# so InsertXML inserts the new paragraphs without disturbing any existing
# content.

$d = $word.ActiveDocument

$newParagraphsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="0000FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>SELECT</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:tab/>
        <w:t xml:space="preserve"> assunto</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="808080"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>ano</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="808080"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="FF00FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>count</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="808080"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>assunto</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="808080"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="0000FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>as</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> Quantidade</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="0000FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>FROM</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:tab/>
        <w:t>atendimentos</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="0000FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>GROUP</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="0000FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>BY</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:tab/>
        <w:t>assunto</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="808080"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>ano</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="0000FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>ORDER</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="0000FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>BY</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="FF00FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>count</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="808080"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>assunto</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="808080"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="0000FF"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>desc</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="808080"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>ano</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$insertionPoint.InsertXML($newParagraphsXml)
